$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $range = $ws.Range($cell)
    $origStyle = $range.Style
    $range.Value = "'" + $value
    $range.Style = $origStyle
}

Set-TextValue 'D2' '25.920.26'
Set-TextValue 'E2' '  +0.38%  '

Set-TextValue 'D3' '1.642.60'
Set-TextValue 'E3' '  +0.65%  '

Set-TextValue 'E4' '  +0.52%  '

Set-TextValue 'D5' '215.35'
Set-TextValue 'E5' '  +0.32%  '

Set-TextValue 'D6' '0.5088'
Set-TextValue 'E6' '  +1.66%  '

Set-TextValue 'D7' '1.005'
Set-TextValue 'E7' '  +0.37%  '

Set-TextValue 'D8' '0.2569'
Set-TextValue 'E8' '  +0.47%  '

Set-TextValue 'D9' '0.06396'
Set-TextValue 'E9' '  +0.74%  '

Set-TextValue 'D10' '19.56'
Set-TextValue 'E10' '  -0.20%  '

Set-TextValue 'D11' '0.07781'
Set-TextValue 'E11' '  +1.21%  '

Set-TextValue 'D12' '4.306'
Set-TextValue 'E12' '  +1.37%  '

Set-TextValue 'D13' '1.653.56'
Set-TextValue 'E13' '  +0.11%  '

Set-TextValue 'D14' '0.5459'
Set-TextValue 'E14' '  +0.88%  '

Set-TextValue 'D15' '0.0₅7858'
Set-TextValue 'E15' '  -0.17%  '

Set-TextValue 'D16' '64.62'
Set-TextValue 'E16' '  +1.75%  '

Set-TextValue 'D17' '25.989.07'
Set-TextValue 'E17' '  +0.67%  '

Set-TextValue 'D18' '1.005'
Set-TextValue 'E18' '  +0.38%  '

Set-TextValue 'D19' '197.93'
Set-TextValue 'E19' '  -1.17%  '

Set-TextValue 'D20' '4.442'
Set-TextValue 'E20' '  +2.95%  '

Set-TextValue 'D21' '9.967'
Set-TextValue 'E21' '  +1.11%  '

Set-TextValue 'D22' '6.033'
Set-TextValue 'E22' '  +1.81%  '

Set-TextValue 'D23' '1.007'
Set-TextValue 'E23' '  +0.47%  '

Set-TextValue 'D24' '1.877'
Set-TextValue 'E24' '  -2.53%  '

Set-TextValue 'D25' '140.57'
Set-TextValue 'E25' '  -0.34%  '

Set-TextValue 'D26' '0.1147'
Set-TextValue 'E26' '  +1.45%  '

Set-TextValue 'D27' '6.908'
Set-TextValue 'E27' '  +3.47%  '

Set-TextValue 'D28' '15.73'
Set-TextValue 'E28' '  +0.92%  '

Set-TextValue 'D29' '1.240'
Set-TextValue 'E29' '  +0.21%  '

Set-TextValue 'D30' '0.05023'
Set-TextValue 'E30' '  +1.02%  '

Set-TextValue 'D31' '3.261'
Set-TextValue 'E31' '  -0.03%  '

Set-TextValue 'D32' '3.189'
Set-TextValue 'E32' '  +0.23%  '

Set-TextValue 'D33' '1.541'
Set-TextValue 'E33' '  +0.57%  '

Set-TextValue 'D34' '2.363'
Set-TextValue 'E34' '  -0.12%  '

Set-TextValue 'D35' '0.8943'
Set-TextValue 'E35' '  +0.82%  '

Set-TextValue 'D36' '2.588'
Set-TextValue 'E36' '  -1.29%  '

Set-TextValue 'D37' '1.133.16'
Set-TextValue 'E37' '  -2.66%  '

Set-TextValue 'D38' '0.5513'
Set-TextValue 'E38' '  -0.57%  '

Set-TextValue 'D39' '0.01559'
Set-TextValue 'E39' '  +0.26%  '

Set-TextValue 'B40' 'PaxDollar'
Set-TextValue 'C40' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D40' '1.005'
Set-TextValue 'E40' '  +0.43%  '

Set-TextValue 'B41' 'mCoin'
Set-TextValue 'C41' 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
Set-TextValue 'D41' '2.552'
Set-TextValue 'E41' '  -0.25%  '

Set-TextValue 'D42' '5.634'
Set-TextValue 'E42' '  -0.35%  '

Set-TextValue 'E43' '  +10.20%  '

Set-TextValue 'D44' '0.8167'
Set-TextValue 'E44' '  +1.93%  '

Set-TextValue 'D45' '99.69'
Set-TextValue 'E45' '  +0.53%  '

Set-TextValue 'D46' '1.782.35'
Set-TextValue 'E46' '  +0.76%  '

Set-TextValue 'E47' '  +0.46%  '

Set-TextValue 'E48' '  +0.47%  '

Set-TextValue 'D49' '54.95'
Set-TextValue 'E49' '  +0.93%  '

Set-TextValue 'D50' '0.05087'
Set-TextValue 'E50' '  +0.41%  '

Set-TextValue 'E51' '  +0.39%  '
